# "Entrega Ejemplo - laboratorio 4"
# Replaces the placeholder timing figures in the two result tables
# ("Tamano de la muestra (ARRAYLIST)" -> A1:D11, and
#  "Tamano de la muestra (LINKED_LIST)" -> A14:D24) on sheet "Datos Lab4"
# with the real measured millisecond values, as plain numbers (no more
# the old +N / running-sum demo formulas), and leaves the not-yet-run
# larger sample sizes blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos Lab4")

# ---- Table 1 : ARRAYLIST (rows 2-11) ----------------------------------
# A column (sample sizes) is untouched.

$ws.Range("B2").Value = 563.12
$ws.Range("C2").Value = 645.57
$ws.Range("D2").Value = 36.7

$ws.Range("B3").Value = 2293.01
$ws.Range("C3").Value = 2593.07
$ws.Range("D3").Value = 77.56

$ws.Range("B4").Value = 9248.66
$ws.Range("C4").Value = 10749.37
$ws.Range("D4").Value = 171.08

$ws.Range("B5").Value = 38467.51
$ws.Range("C5").Value = 45474.61
$ws.Range("D5").Value = 417.92

$ws.Range("B6").Value = 164240.64
$ws.Range("C6").Value = 190602.47
$ws.Range("D6").Value = 913.15

$ws.Range("B7").Value = 688699.79
$ws.Range("C7").Value = 749434.38
$ws.Range("D7").Value = 2158.07

$ws.Range("B8:C8").ClearContents()
$ws.Range("D8").Value = 5339.5

$ws.Range("B9:C9").ClearContents()
$ws.Range("D9").Value = 12081.23

$ws.Range("B10:C10").ClearContents()
$ws.Range("D10").Value = 31470.36

$ws.Range("B11:D11").ClearContents()

# ---- Table 2 : LINKED_LIST (rows 15-24) --------------------------------

$ws.Range("B15").Value = 44019.46
$ws.Range("C15").Value = 39137.78
$ws.Range("D15").Value = 2196.17

$ws.Range("B16").Value = 362439.99
$ws.Range("C16").Value = 321500.04
$ws.Range("D16").Value = 10308.12

$ws.Range("B17").Value = 2996884.79
$ws.Range("C17").Value = 2615806.89
$ws.Range("D17").Value = 48337.59

$ws.Range("B18:C18").ClearContents()
$ws.Range("D18").Value = 241678.25

$ws.Range("B19:C19").ClearContents()
$ws.Range("D19").Value = 1075527.61

$ws.Range("B20:D20").ClearContents()
$ws.Range("B21:D21").ClearContents()
$ws.Range("B22:D22").ClearContents()
$ws.Range("B23:D23").ClearContents()
$ws.Range("B24:D24").ClearContents()

# ---- Row heights for the two header rows -------------------------------
$ws.Rows.Item(1).RowHeight = 16
$ws.Rows.Item(14).RowHeight = 16

# ---- Selection, matching the cursor position left in the saved file ----
$ws.Range("F13").Select()
